$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.282.06'
$ws.Range('E2').Value = '  +2.55%  '
$ws.Range('D3').Value = '1.871.37'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '338.10'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4699'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3929'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.91%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.09'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.76'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').Value = '1.884.60'
$ws.Range('E13').Value = '  +2.73%  '
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.277'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.12'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001043'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06581'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.66'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.65%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').Value = '28.294.17'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.291'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').Value = '2.113.03'
$ws.Range('E26').Value = '  +2.68%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '159.37'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.86'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.157'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.498'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.81%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '120.03'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9795'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09487'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.590'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.381'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.354'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02268'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06090'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.442'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.177'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5965'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1877'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.302'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.74%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5613'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.16'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.968'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06900'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '110.75'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.016'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +13.22%  '
